$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.740.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.427.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.539"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.415.29"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("E10").Value = "  -4.01%  "
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.813.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.19%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000166"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.833.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +13.44%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.335.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.08%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "559.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -11.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.530.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  -4.71%  "
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "151.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -2.05%  "
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0287"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.588"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0898"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.29%  "
